$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.475.52'
$ws.Range('E2').Value = '  +7.48%  '

$ws.Range('D3').Value = '2.383.61'
$ws.Range('E3').Value = '  +4.88%  '

$ws.Range('E4').Value = '  -0.73%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '112.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +9.11%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '317.88'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.84%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.638'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.91%  '

$ws.Range('E8').Value = '  -0.30%  '

$ws.Range('E9').Value = '  +5.83%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.28'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +10.35%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0931'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.13%  '

$ws.Range('E12').Value = '  +6.56%  '

$ws.Range('E13').Value = '  +5.20%  '

$ws.Range('E14').Value = '  +1.53%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.82'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.93%  '

$ws.Range('D16').Value = '2.744.12'
$ws.Range('E16').Value = '  +4.90%  '

$ws.Range('D17').Value = '2.385.47'
$ws.Range('E17').Value = '  +4.13%  '

$ws.Range('D18').Value = '45.421.38'
$ws.Range('E18').Value = '  +7.42%  '

$ws.Range('E19').Value = '  +6.85%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000109'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.21%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.52%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '75.11'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.61%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.55'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.87%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '270.97'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.78%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.34'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +8.66%  '

$ws.Range('E26').Value = '  -0.62%  '

$ws.Range('E27').Value = '  +7.14%  '

$ws.Range('E28').Value = '  +10.03%  '

$ws.Range('E29').Value = '  +0.27%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.43%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '38.68'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +9.52%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0947'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +11.94%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '169.92'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.97%  '

$ws.Range('E34').Value = '  +17.63%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.134'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.79%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.90'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.91%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.118'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.18%  '

$ws.Range('E38').Value = '  +13.33%  '

$ws.Range('E39').Value = '  +6.05%  '

$ws.Range('E40').Value = '  +8.07%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.74'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +13.42%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '105.50'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +7.11%  '

$ws.Range('E43').Value = '  +8.15%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.49'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +14.64%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '71.28'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.99%  '

$ws.Range('E46').Value = '  +0.18%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '118.61'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +8.76%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.86'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +15.13%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.65'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +22.04%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +8.82%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.21'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.92%  '
